# "Updated for launcPage test case"
#
# - LaunchPageScript becomes the active sheet/tab (was LaunchPage).
# - LaunchPageScript's 4-column test-script table is rewritten: new header
#   wording + bold header row, and 3 fresh scripted-test rows (in addition
#   to the existing row 2, whose copy/description also changes) with
#   wrapped text.
# - Column widths / page setup tidied to fit the longer wrapped text.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("LaunchPage")
$ws4 = $wb.Worksheets.Item("LaunchPageScript")

# ---- Header row (row 1): new captions, bold ----------------------------
$ws4.Range("A1").Value = "Description"
$ws4.Range("B1").Value = "Expected Result"
$ws4.Range("C1").Value = "Actual Pass Result"
$ws4.Range("D1").Value = "Actual Fail Result"
$ws4.Range("A1:D1").Font.Bold = $true

# ---- Row 2: existing "Open Swiggy HomePage" case -> "Launch URL" case --
$ws4.Range("A2").Value = "Launch URL https://www.swiggy.com/"
$ws4.Range("B2").Value = "Swiggy Home Page should open."
$ws4.Range("C2").Value = "Swiggy Home Page opened successfully"
$ws4.Range("D2").Value = "Problem loading swiggy home page."
$ws4.Range("B2").WrapText = $true
$ws4.Rows.Item(2).RowHeight = 45

# ---- Row 3: login / signup links case -----------------------------------
$ws4.Range("A3").Value = "Verify that the login and signup links are present`n on the homepage."
$ws4.Range("B3").Value = "Login and signp links should be present."
$ws4.Range("C3").Value = "Login and signup links are present."
$ws4.Range("D3").Value = "Problem finding the login and signup links on the homepage."
$ws4.Range("A3:D3").WrapText = $true
$ws4.Rows.Item(3).RowHeight = 30

# ---- Row 4: delivery-location input box case -----------------------------
$ws4.Range("A4").Value = "Verify that the input box for entering delivery `nlocation is present along with Locate Me and `nFind Food buttons. "
$ws4.Range("B4").Value = "Input box for entering delivery location should be present along with Locate Me and Find Food buttons."
$ws4.Range("C4").Value = "Input box for entering delivery location is present along with Locate Me and Find Food buttons"
$ws4.Range("D4").Value = "Problem finding the Input box along with Locate Me and Find Food buttons."
$ws4.Range("A4:D4").WrapText = $true
$ws4.Rows.Item(4).RowHeight = 60

# ---- Row 5: text-entry case ------------------------------------------------
$ws4.Range("A5").Value = "Verify that text can be entered in the inputbox fordelivery location."
$ws4.Range("B5").Value = "Text should be entered in the input box for delivery location."
$ws4.Range("C5").Value = "Text  entered in the input box for delivery location."
$ws4.Range("D5").Value = "Problem entering text in the input box for delivery location."
$ws4.Range("A5:D5").WrapText = $true
$ws4.Rows.Item(5).RowHeight = 30

# ---- Column widths (sized to fit the new, longer wrapped text) ----------
$ws4.Columns.Item(1).ColumnWidth = 44.166666666666664
$ws4.Columns.Item(2).ColumnWidth = 29.5
$ws4.Columns.Item(3).ColumnWidth = 35.83333333333333
$ws4.Columns.Item(4).ColumnWidth = 32.666666666666664

# ---- Page setup ------------------------------------------------------------
$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1

# ---- Selection / active sheet ---------------------------------------------
$ws4.Range("D5").Select()
$ws4.Activate()

$wb.Save()
